# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Latest Handoff Datetime" on each localized-language sheet for the row
# corresponding to file 796cd5df-4a86-45f3-a5d2-37fbc8860813 (the 6th data
# row / sheet row 7), reflecting a freshly generated handoff package.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G7 - "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-09-03 08:45:07"

# zh-cn!H7 - "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-09-03 08:44:58"

# de-de!H7 - "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-09-03 08:45:07"
